# Auto-generated edit script: applies scheduled-runner market-data refresh
# to the Gungnir_Profits workbook's per-class leve-crafting sheets.
# For each changed cell we set the new numeric value directly via the
# Excel COM object model (Range.Value), matching the OOXML diff exactly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 76.666664
$ws.Range("I12").Value = 80
$ws.Range("J12").Value = 70
$ws.Range("K12").Value = 80
$ws.Range("L12").Value = 70
$ws.Range("M12").Value = 90
$ws.Range("N12").Value = -410
$ws.Range("H137").Value = 1157.2307
$ws.Range("I137").Value = 1033.3226
$ws.Range("J137").Value = 1637.375
$ws.Range("K137").Value = 3099.9678
$ws.Range("L137").Value = 4912.125
$ws.Range("M137").Value = -549.9677999999999
$ws.Range("N137").Value = -10012.125
$ws.Range("H138").Value = 1666.2693
$ws.Range("I138").Value = 811.84375
$ws.Range("J138").Value = 3033.35
$ws.Range("K138").Value = 2435.53125
$ws.Range("L138").Value = 9100.049999999999
$ws.Range("M138").Value = 2704.46875
$ws.Range("N138").Value = -19380.05
$ws.Range("H141").Value = 2682.5715
$ws.Range("I141").Value = 1758.8462
$ws.Range("J141").Value = 5351.1113
$ws.Range("K141").Value = 5276.5386
$ws.Range("L141").Value = 16053.3339
$ws.Range("M141").Value = -96.53859999999986
$ws.Range("N141").Value = -26413.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 5750
$ws.Range("J6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("N6").Value = -10346
$ws.Range("H32").Value = 797.14
$ws.Range("I32").Value = 760.14435
$ws.Range("J32").Value = 1993.3334
$ws.Range("K32").Value = 760.14435
$ws.Range("L32").Value = 1993.3334
$ws.Range("M32").Value = -473.14435
$ws.Range("N32").Value = -2567.3334
$ws.Range("H45").Value = 38126
$ws.Range("I45").Value = 53607.316
$ws.Range("J45").Value = 1357.875
$ws.Range("K45").Value = 53607.316
$ws.Range("L45").Value = 1357.875
$ws.Range("M45").Value = -53230.316
$ws.Range("N45").Value = -2111.875
$ws.Range("H61").Value = 1181.12
$ws.Range("I61").Value = 1172
$ws.Range("J61").Value = 1400
$ws.Range("K61").Value = 1172
$ws.Range("L61").Value = 1400
$ws.Range("M61").Value = -960
$ws.Range("N61").Value = -1824
$ws.Range("H121").Value = 22400
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 22400
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 22400
$ws.Range("N121").Value = -25894
$ws.Range("H122").Value = 840
$ws.Range("I122").Value = 840
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2520
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -70
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H124").Value = 34559.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 34559.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 34559.5
$ws.Range("N124").Value = -44379.5
$ws.Range("H125").Value = 61983
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 61983
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 61983
$ws.Range("N125").Value = -71823
$ws.Range("H126").Value = 5610
$ws.Range("I126").Value = 5610
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 16830
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -14360
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 49800
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 49800
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 49800
$ws.Range("N128").Value = -59760
$ws.Range("H129").Value = 25040.375
$ws.Range("I129").Value = 10909
$ws.Range("J129").Value = 48592.668
$ws.Range("K129").Value = 10909
$ws.Range("L129").Value = 48592.668
$ws.Range("M129").Value = -5909
$ws.Range("N129").Value = -58592.668
$ws.Range("H130").Value = 41000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 41000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 41000
$ws.Range("N130").Value = -51040
$ws.Range("H131").Value = 80357.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 80357.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 80357.5
$ws.Range("N131").Value = -90437.5
$ws.Range("H132").Value = 1402050.5
$ws.Range("I132").Value = 1110.5714
$ws.Range("J132").Value = 4203930.5
$ws.Range("K132").Value = 3331.7142
$ws.Range("L132").Value = 12611791.5
$ws.Range("M132").Value = -801.7142000000003
$ws.Range("N132").Value = -12616851.5
$ws.Range("H133").Value = 37000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 37000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 37000
$ws.Range("N133").Value = -42060
$ws.Range("H134").Value = 30413.166
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 30413.166
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 30413.166
$ws.Range("N134").Value = -40553.166
$ws.Range("H135").Value = 40489.855
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 40489.855
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 40489.855
$ws.Range("N135").Value = -50629.855
$ws.Range("H136").Value = 1181.12
$ws.Range("I136").Value = 1172
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 3516
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -966
$ws.Range("N136").Value = -9300
$ws.Range("H137").Value = 46972
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 46972
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 46972
$ws.Range("N137").Value = -57172
$ws.Range("H138").Value = 47493
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 47493
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 47493
$ws.Range("N138").Value = -57773
$ws.Range("H139").Value = 42517.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 42517.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 42517.5
$ws.Range("N139").Value = -52797.5
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 45455390
$ws.Range("I107").Value = 50000630
$ws.Range("K107").Value = 50000630
$ws.Range("M107").Value = -49998710

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1332.6852
$ws.Range("I31").Value = 962.95
$ws.Range("J31").Value = 1550.1765
$ws.Range("K31").Value = 962.95
$ws.Range("L31").Value = 1550.1765
$ws.Range("M31").Value = -667.95
$ws.Range("N31").Value = -2140.1765
$ws.Range("H34").Value = 1332.6852
$ws.Range("I34").Value = 962.95
$ws.Range("J34").Value = 1550.1765
$ws.Range("K34").Value = 962.95
$ws.Range("L34").Value = 1550.1765
$ws.Range("M34").Value = -760.95
$ws.Range("N34").Value = -1954.1765
$ws.Range("H58").Value = 19231394
$ws.Range("I58").Value = 27778302
$ws.Range("J58").Value = 848.625
$ws.Range("K58").Value = 27778302
$ws.Range("L58").Value = 848.625
$ws.Range("M58").Value = -27778099
$ws.Range("N58").Value = -1254.625
$ws.Range("H132").Value = 7408635.5
$ws.Range("I132").Value = 873.0645
$ws.Range("J132").Value = 23811538
$ws.Range("K132").Value = 2619.1935
$ws.Range("L132").Value = 71434614
$ws.Range("M132").Value = -89.19349999999986
$ws.Range("N132").Value = -71439674
$ws.Range("H136").Value = 19231394
$ws.Range("I136").Value = 27778302
$ws.Range("J136").Value = 848.625
$ws.Range("K136").Value = 83334906
$ws.Range("L136").Value = 2545.875
$ws.Range("M136").Value = -83332356
$ws.Range("N136").Value = -7645.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16669634
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 33335268
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 33335268
$ws.Range("M80").Value = -3002
$ws.Range("N80").Value = -33337264
$ws.Range("H83").Value = 16669634
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 33335268
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 166676340
$ws.Range("M83").Value = -15008
$ws.Range("N83").Value = -166686324
$ws.Range("H132").Value = 14679.883
$ws.Range("I132").Value = 13610.25
$ws.Range("J132").Value = 15630.667
$ws.Range("K132").Value = 40830.75
$ws.Range("L132").Value = 46892.001
$ws.Range("M132").Value = -38300.75
$ws.Range("N132").Value = -51952.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 1461.2
$ws.Range("I12").Value = 434
$ws.Range("J12").Value = 3002
$ws.Range("K12").Value = 434
$ws.Range("L12").Value = 3002
$ws.Range("M12").Value = -264
$ws.Range("N12").Value = -3342
$ws.Range("H132").Value = 17246648
$ws.Range("I132").Value = 27779240
$ws.Range("K132").Value = 83337720
$ws.Range("M132").Value = -83335190
$ws.Range("H136").Value = 39412056
$ws.Range("I136").Value = 7522825.5
$ws.Range("J136").Value = 100001590
$ws.Range("K136").Value = 22568476.5
$ws.Range("L136").Value = 300004770
$ws.Range("M136").Value = -22565926.5
$ws.Range("N136").Value = -300009870

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 20469.703
$ws.Range("I122").Value = 32635.188
$ws.Range("J122").Value = 2774.4546
$ws.Range("K122").Value = 97905.564
$ws.Range("L122").Value = 8323.363799999999
$ws.Range("M122").Value = -95455.564
$ws.Range("N122").Value = -13223.3638
$ws.Range("H132").Value = 22924.076
$ws.Range("I132").Value = 33551.03
$ws.Range("J132").Value = 6730.619
$ws.Range("K132").Value = 100653.09
$ws.Range("L132").Value = 20191.857
$ws.Range("M132").Value = -98123.09
$ws.Range("N132").Value = -25251.857

Write-Output "Applied 270 cell updates across 7 sheets."
